$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates (crypto prices / 1h volume changes)
# Row 2
$ws.Range("D2").Value = '48.373.34'
$ws.Range("E2").Value = '  +2.43%  '

# Row 3
$ws.Range("D3").Value = '2.522.92'
$ws.Range("E3").Value = '  +1.44%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '109.78'
$ws.Range("E5").Value = '  +1.92%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '322.39'
$ws.Range("E6").Value = '  +0.44%  '

# Row 7
$ws.Range("E7").Value = '  +2.18%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("E9").Value = '  +3.76%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.56'
$ws.Range("E10").Value = '  +5.43%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.42'
$ws.Range("E11").Value = '  +12.13%  '

# Row 12
$ws.Range("E12").Value = '  +1.73%  '

# Row 13
$ws.Range("E13").Value = '  +1.10%  '

# Row 14
$ws.Range("E14").Value = '  +2.56%  '

# Row 15
$ws.Range("D15").Value = '2.920.23'
$ws.Range("E15").Value = '  +1.62%  '

# Row 16
$ws.Range("D16").Value = '2.517.68'
$ws.Range("E16").Value = '  +1.13%  '

# Row 17
$ws.Range("E17").Value = '  +1.09%  '

# Row 18
$ws.Range("D18").Value = '48.195.48'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.33'
$ws.Range("E19").Value = '  +4.99%  '

# Row 20
$ws.Range("E20").Value = '  +0.39%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0950'
$ws.Range("E21").Value = '  +2.18%  '

# Row 22
$ws.Range("E22").Value = '  -0.96%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.05'
$ws.Range("E23").Value = '  +2.60%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '273.37'
$ws.Range("E24").Value = '  +11.55%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.59'
$ws.Range("E25").Value = '  +0.93%  '

# Row 26
$ws.Range("E26").Value = '  +0.00%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.01'
$ws.Range("E27").Value = '  +1.42%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.38'
$ws.Range("E28").Value = '  +4.57%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.13'
$ws.Range("E29").Value = '  +1.41%  '

# Row 30
$ws.Range("E30").Value = '  +7.19%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.69'
$ws.Range("E31").Value = '  +4.03%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.66'
$ws.Range("E32").Value = '  +0.42%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.74'
$ws.Range("E33").Value = '  -2.78%  '

# Row 34
$ws.Range("E34").Value = '  +1.37%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.01'
$ws.Range("E35").Value = '  +0.10%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0786'
$ws.Range("E36").Value = '  +1.24%  '

# Row 37
$ws.Range("E37").Value = '  +1.40%  '

# Row 38
$ws.Range("E38").Value = '  +1.67%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.01'
$ws.Range("E39").Value = '  +3.17%  '

# Row 40
$ws.Range("E40").Value = '  +0.76%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '122.92'
$ws.Range("E41").Value = '  +3.39%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.23'
$ws.Range("E42").Value = '  -1.85%  '

# Row 43
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.22'
$ws.Range("E43").Value = '  +0.11%  '

# Row 44
$ws.Range("E44").Value = '  +2.41%  '

# Row 45
$ws.Range("D45").Value = '2.028.86'
$ws.Range("E45").Value = '  +2.33%  '

# Row 46
$ws.Range("E46").Value = '  +5.04%  '

# Row 47
$ws.Range("E47").Value = '  +7.69%  '

# Row 48
$ws.Range("E48").Value = '  +2.89%  '

# Row 49
$ws.Range("E49").Value = '  +1.11%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.21'
$ws.Range("E50").Value = '  +1.91%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.83'
$ws.Range("E51").Value = '  +4.11%  '
